$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new server-data row (row 2).
# Header mapping (row 1): A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "LoginServer_1"
$ws.Range("B2").Value = "000106001"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "LoginServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 6001

# Data validation list on F2:F1048576 is trimmed to start from F3 (F2 now holds real data).
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Selection moves from A2:H7 to G5.
$ws.Range("G5").Select()
